# Updated cryptos list on Thu Oct 17 03:55:58 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) columns of the
# cryptos listing sheet with newly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => { D = new price text (optional), E = new volume text }
$updates = @(
    @{ Row = 2;  D = "67.624.47";   E = "  +0.70%  " },
    @{ Row = 3;  D = "2.642.04";    E = "  +0.81%  " },
    @{ Row = 5;  D = "604.51";      E = "  +1.67%  " },
    @{ Row = 6;  D = "154.66";      E = "  +0.18%  " },
    @{ Row = 7;  E = "  +0.00%  " },
    @{ Row = 8;  D = "0.549";       E = "  +1.07%  " },
    @{ Row = 9;  D = "2.641.57";    E = "  +0.81%  " },
    @{ Row = 10; E = "  +7.59%  " },
    @{ Row = 11; E = "  +0.52%  " },
    @{ Row = 12; D = "5.23";        E = "  +0.06%  " },
    @{ Row = 13; D = "0.353";       E = "  -0.84%  " },
    @{ Row = 14; D = "28.07";       E = "  +1.07%  " },
    @{ Row = 15; D = "3.121.46";    E = "  +0.83%  " },
    @{ Row = 16; D = "0.0000185";   E = "  +1.35%  " },
    @{ Row = 17; D = "67.533.84";   E = "  +0.78%  " },
    @{ Row = 18; D = "2.638.98";    E = "  +0.78%  " },
    @{ Row = 19; D = "11.34";       E = "  +0.21%  " },
    @{ Row = 20; D = "365.84";      E = "  +0.89%  " },
    @{ Row = 21; D = "7.64";        E = "  -4.23%  " },
    @{ Row = 22; E = "  -0.45%  " },
    @{ Row = 23; E = "  +7.78%  " },
    @{ Row = 24; E = "  +0.13%  " },
    @{ Row = 25; D = "10.04";       E = "  -2.36%  " },
    @{ Row = 26; D = "66.14";       E = "  -7.98%  " },
    @{ Row = 28; D = "2.758.91";    E = "  +0.30%  " },
    @{ Row = 29; D = "584.02";      E = "  -7.43%  " },
    @{ Row = 30; E = "  +0.30%  " },
    @{ Row = 31; D = "1.43";        E = "  -2.80%  " },
    @{ Row = 32; E = "  -1.22%  " },
    @{ Row = 33; E = "  -0.20%  " },
    @{ Row = 34; E = "  -2.28%  " },
    @{ Row = 36; E = "  -1.87%  " },
    @{ Row = 37; D = "4.98";        E = "  -0.50%  " },
    @{ Row = 38; D = "158.05";      E = "  +1.95%  " },
    @{ Row = 39; D = "19.52";       E = "  +0.07%  " },
    @{ Row = 40; E = "  +0.50%  " },
    @{ Row = 41; E = "  -3.48%  " },
    @{ Row = 42; E = "  +0.14%  " },
    @{ Row = 43; D = "2.64";        E = "  +0.86%  " },
    @{ Row = 44; D = "41.17";       E = "  -0.47%  " },
    @{ Row = 45; E = "  -0.03%  " },
    @{ Row = 46; E = "  -0.73%  " },
    @{ Row = 47; D = "156.71";      E = "  +0.41%  " },
    @{ Row = 48; E = "  -3.54%  " },
    @{ Row = 49; E = "  -1.02%  " },
    @{ Row = 51; D = "0.629";       E = "  +0.46%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        # Several "price" strings (e.g. "604.51") are valid numeric literals.
        # Force the cell to Text so Excel stores the exact original string
        # instead of silently converting it to a number, then restore the
        # cell's previous (default/general) number format so no stray
        # style is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
